# Update "想去人数" (column F) counts across the four sheets to reflect
# the refreshed scrape (gh-pages data regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1803
$ws1.Range("F4").Value  = 473
$ws1.Range("F7").Value  = 636
$ws1.Range("F8").Value  = 353
$ws1.Range("F9").Value  = 1764
$ws1.Range("F10").Value = 382
$ws1.Range("F11").Value = 1437
$ws1.Range("F12").Value = 823
$ws1.Range("F13").Value = 348
$ws1.Range("F14").Value = 694
$ws1.Range("F15").Value = 12918
$ws1.Range("F16").Value = 12881
$ws1.Range("F17").Value = 964
$ws1.Range("F18").Value = 749
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 528
$ws1.Range("F21").Value = 55
$ws1.Range("F22").Value = 591
$ws1.Range("F23").Value = 2021
$ws1.Range("F24").Value = 38
$ws1.Range("F25").Value = 15
$ws1.Range("F27").Value = 17
$ws1.Range("F28").Value = 102
$ws1.Range("F30").Value = 690

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 87
$ws2.Range("F10").Value = 84

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 174

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 174
$ws4.Range("F5").Value  = 1803
$ws4.Range("F6").Value  = 473
$ws4.Range("F11").Value = 636
$ws4.Range("F13").Value = 353
$ws4.Range("F14").Value = 1764
$ws4.Range("F15").Value = 382
$ws4.Range("F16").Value = 1437
$ws4.Range("F17").Value = 823
$ws4.Range("F18").Value = 348
$ws4.Range("F19").Value = 87
$ws4.Range("F20").Value = 694
$ws4.Range("F21").Value = 12918
$ws4.Range("F22").Value = 12881
$ws4.Range("F23").Value = 964
$ws4.Range("F24").Value = 749
$ws4.Range("F25").Value = 12
$ws4.Range("F26").Value = 528
$ws4.Range("F27").Value = 55
$ws4.Range("F28").Value = 591
$ws4.Range("F31").Value = 2021
$ws4.Range("F32").Value = 38
$ws4.Range("F33").Value = 15
$ws4.Range("F36").Value = 17
$ws4.Range("F38").Value = 102
$ws4.Range("F40").Value = 690
$ws4.Range("F41").Value = 84
